$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-changed date for each record.
# Every populated data row (2 through 106) currently shows 45175
# (2023-09-06); bump it to 45177 (2023-09-08), leaving everything else
# (formatting, other columns) untouched.
for ($row = 2; $row -le 106; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45175) {
        $cell.Value = 45177
    }
}
